$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: "dimenzija za m=4" (columns A,B) ---
$sheet1Data = @(
    @(41, 17),
    @(42, 17),
    @(43, 18),
    @(44, 18),
    @(45, 18),
    @(46, 19),
    @(47, 19),
    @(48, 20),
    @(49, 20),
    @(50, 20),
    @(51, 21),
    @(52, 21),
    @(53, 22),
    @(54, 22),
    @(55, 22),
    @(56, 23),
    @(57, 23),
    @(58, 24),
    @(59, 24),
    @(60, 24),
    @(61, 25),
    @(62, 25),
    @(63, 26),
    @(64, 26),
    @(65, 26),
    @(66, 27),
    @(67, 27),
    @(68, 28),
    @(69, 28),
    @(70, 28),
    @(71, 29),
    @(72, 29),
    @(73, 30),
    @(74, 30),
    @(75, 30),
    @(76, 31),
    @(77, 31)
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $r = 36 + $i
    $row = $sheet1Data[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
}

# --- Sheet 2: "dimenzija za m=5-10" (columns A-G) ---
$sheet2Data = @(
    @(41, 15, 16, 16, 16, 15, 15),
    @(42, 14, 16, 16, 14, 17, 16),
    @(43, 16, 17, 17, 17, 17, 16),
    @(44, 16, 16, 17, 17, 16, 16),
    @(45, 15, 18, 17, 15, 18, 18),
    @(46, 17, 18, 18, 18, 18, 17),
    @(47, 17, 18, 18, 18, 18, 18),
    @(48, 16, 18, 18, 16, 18, 18),
    @(49, 18, 19, 19, 19, 19, 18),
    @(50, 18, 19, 19, 19, 20, 20),
    @(51, 17, 19, 20, 17, 18, 19),
    @(52, 19, 20, 20, 20, 21, 19),
    @(53, 19, 20, 20, 20, 20, 20),
    @(54, 18, 21, 21, 18, 20, 20),
    @(55, 20, 20, 21, 21, 21, 20),
    @(56, 20, 22, 21, 21, 21, 21),
    @(57, 19, 22, 22, 19, 22, 21),
    @(58, 21, 22, 22, 22, 21, 22),
    @(59, 21, 22, 23, 22, 23, 22),
    @(60, 20, 23, 23, 20, 23, 22),
    @(61, 22, 23, 23, 23, 22, 23),
    @(62, 22, 23, 24, 23, 24, 23),
    @(63, 21, 24, 24, 21, 24, 23),
    @(64, 23, 24, 24, 24, 24, 24),
    @(65, 23, 25, 25, 24, 24, 24),
    @(66, 22, 24, 25, 22, 25, 24),
    @(67, 24, 26, 26, 25, 26, 25),
    @(68, 24, 26, 26, 25, 24, 25),
    @(69, 23, 26, 26, 23, 27, 26),
    @(70, 25, 26, 27, 26, 26, 26),
    @(71, 25, 27, 27, 26, 26, 26),
    @(72, 24, 27, 27, 24, 27, 27),
    @(73, 26, 27, 28, 27, 27, 27),
    @(74, 26, 28, 28, 27, 28, 27),
    @(75, 25, 28, 29, 25, 27, 28),
    @(76, 27, 29, 29, 28, 29, 28),
    @(77, 27, 28, 29, 28, 29, 28)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $r = 36 + $i
    $row = $sheet2Data[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws2.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# --- Copy the bordered/bold/centered style from column A header rows down to the new rows ---
$ws1.Range("A2").Copy()
$ws1.Range("A36:A72").PasteSpecial(-4122)
$ws2.Range("A2").Copy()
$ws2.Range("A36:A72").PasteSpecial(-4122)

$excel.CutCopyMode = $false
